$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25: new row, cloned from the former row 24's content (A,B,C,E,F,G,H,I,J,N,Q,R
#     unchanged; D/K/L/M/O/P are the values the old row 24 carried) ---
$ws.Cells.Item(25, 1).Value = 11
$ws.Cells.Item(25, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(25, 3).Value = "Bíobío"
$ws.Cells.Item(25, 4).Value = 44272
$ws.Range("D25").NumberFormat = $ws.Range("D24").NumberFormat()
$ws.Cells.Item(25, 5).Value = 8
$ws.Cells.Item(25, 6).Value = 100112030
$ws.Cells.Item(25, 7).Value = "Poroto granado"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 100
$ws.Cells.Item(25, 11).Value = 22000
$ws.Cells.Item(25, 12).Value = 24000
$ws.Cells.Item(25, 13).Value = 23000
$ws.Cells.Item(25, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(25, 15).Value = "Región del Maule"
$ws.Cells.Item(25, 16).Value = 920
$ws.Cells.Item(25, 17).Value = 25
$ws.Cells.Item(25, 18).Value = "Hortaliza"

# --- Row 24: date moves to what used to be row 23's later date, and the
#     price/origin figures become what used to be row 22's figures ---
$ws.Range("D24").Value = 44313
$ws.Range("K24").Value = 30000
$ws.Range("L24").Value = 32000
$ws.Range("M24").Value = 31000
$ws.Range("O24").Value = "Región Metropolitana"
$ws.Range("P24").Value = 1240

# --- Row 23: only the date changes ---
$ws.Range("D23").Value = 44194

# --- Row 22: new sample week (date, volume, prices, origin) ---
$ws.Range("D22").Value = 44568
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 25000
$ws.Range("L22").Value = 26000
$ws.Range("M22").Value = 25500
$ws.Range("O22").Value = "Región de O'Higgins"
$ws.Range("P22").Value = 1020
